# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# The only substantive content change is that the "Valor Mora" values in
# rows 16 and 24 (column F) were swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 36000
$ws.Range("F24").Value = 72000
